$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vscs")

# ---------------------------------------------------------------------------
# Step 1: Capture the existing comments (text) for every row from 34 to 71
# (the region that will be shifted / rewritten) and remove them. They will be
# re-created after the rows have moved, at their correct new location.
# ---------------------------------------------------------------------------
$commentTexts = @{}
for ($r = 34; $r -le 71; $r++) {
    $c = $ws.Range("A" + $r).Comment
    if ($c -ne $null) {
        $commentTexts[$r] = $c.Text()
        $c.Delete()
    }
}

# ---------------------------------------------------------------------------
# Step 2: Insert 4 new rows right after row 34 (i.e. rows 35-38 become new,
# blank rows and everything that used to live at row 35+ now lives 4 rows
# further down). Excel automatically shifts merged cells and data
# validations that live below the insertion point.
# ---------------------------------------------------------------------------
$ws.Rows("35:38").Insert()

# ---------------------------------------------------------------------------
# Step 3: Re-create the comments that used to live on rows 34-71 on their new
# row (old row 34 is handled separately below, since its text changes).
# ---------------------------------------------------------------------------
foreach ($r in 35..71) {
    if ($commentTexts.ContainsKey($r)) {
        $newRow = $r + 4
        $txt = $commentTexts[$r]
        $ws.Range("A" + $newRow).AddComment($txt) | Out-Null
    }
}

# ---------------------------------------------------------------------------
# Step 4: Fix up the text of row 33 and 34 (renamed in place) and give the
# new rows 35-38 the correct formatting (copied from representative rows
# that already carry the right style) before filling in their text.
# ---------------------------------------------------------------------------
$ws.Range("A33").Value = "BGP Interface"
$ws.Range("A34").Value = "BGP Interface IP Address"
if ($commentTexts.ContainsKey(34)) {
    $ws.Range("A34").AddComment("IP Address for Optional BGP Interface") | Out-Null
}

# Row 35 ("BGP Interface Prefix length") uses the same banded style as row 42
# (formerly row 38, "vCenter VM Folder") -- style s=6/7.
$ws.Range("A42:C42").Copy()
$ws.Range("A35:C35").PasteSpecial(-4122)
$ws.Range("A35").Value = "BGP Interface Prefix length"
$ws.Range("A35").AddComment("Prefix length for the optional BGP interface [default: 24]") | Out-Null

# Row 36 ("BGP Interface VLAN ID") uses the same style as row 34 -- style s=8/9.
$ws.Range("A34:C34").Copy()
$ws.Range("A36:C36").PasteSpecial(-4122)
$ws.Range("A36").Value = "BGP Interface VLAN ID"
$ws.Range("A36").AddComment("VLAN ID for the optional BGP interface [default: 1000]") | Out-Null

# Row 37 is a new section header ("vCenter Parameters"), style s=3, same as
# row 33. Headers only occupy column A (no B/C cells), so copy a single cell.
$ws.Range("A33").Copy()
$ws.Range("A37").PasteSpecial(-4122)
$ws.Range("A37").Value = "vCenter Parameters"
$ws.Range("A37:C37").Merge()

# Row 38 ("vCenter Datacenter Name") uses the same style as row 34 -- s=8/9.
$ws.Range("A34:C34").Copy()
$ws.Range("A38:C38").PasteSpecial(-4122)
$ws.Range("A38").Value = "vCenter Datacenter Name"

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Step 5: Whole-number data validation on the two new numeric fields.
# ---------------------------------------------------------------------------
foreach ($addr in @("B35", "C35", "B36", "C36")) {
    $rng = $ws.Range($addr)
    $rng.Validation.Delete()
    $rng.Validation.Add(1, 1, 1, 0)
    $rng.Validation.ErrorTitle = "Invalid Entry"
    $rng.Validation.ErrorMessage = "Your entry is not an integer, change anyway?"
    $rng.Validation.PromptTitle = "Integer Selection"
    $rng.Validation.InputMessage = "Please provide integer"
    $rng.Validation.IgnoreBlank = 1
    $rng.Validation.ShowInput = 1
    $rng.Validation.ShowError = 1
}

# ---------------------------------------------------------------------------
# Step 6: Keep the sheet's used-range / dimension correct.
# ---------------------------------------------------------------------------
$ws.Range("A1").Select()
